$wb = $excel.ActiveWorkbook

# 1. Insert a new worksheet named "datatypes" before the first existing sheet.
#    Copying an existing sheet (then clearing it) avoids the engine stamping
#    a default sheetFormatPr/baseColWidth that a brand-new blank sheet gets.
$firstSheet = $wb.Worksheets.Item(1)
$firstSheet.Copy($firstSheet)
$ws = $wb.Worksheets.Item(1)
$ws.Name = "datatypes"
$ws.Cells.Clear() | Out-Null

# 2. Populate the new sheet with a variety of data types
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 1.5
$ws.Range("A3").Formula = '=CONCATENATE("a","b")'
$ws.Range("A4").Formula = '=A1>A2'
$ws.Range("A5").Value = "test"
$ws.Range("A6").Value = 42663
$ws.Range("A6").NumberFormat = "mm-dd-yy"

# Autofit column A to match the content
$ws.Columns("A:A").ColumnWidth = 9.75

# Put the selection where it ends up after data entry
$ws.Range("A7").Select() | Out-Null

# 3. issue5 sheet: just move the selection
$issue5 = $wb.Worksheets.Item("issue5")
$issue5.Range("C33").Select() | Out-Null

# 4. issue6 sheet: add a new formatted (but empty) cell A6, resize dimension,
#    and change the selection to A1:A4. Copy/paste the date format from the
#    datatypes sheet so both cells share the same style entry.
$issue6 = $wb.Worksheets.Item("issue6")
$ws.Range("A6").Copy() | Out-Null
$issue6.Range("A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$issue6.Columns("A:A").ColumnWidth = 9.75
$issue6.Range("A1:A4").Select() | Out-Null

# 5. Re-activate the new sheet so it becomes the selected tab, and move the
#    workbook view off the old active tab (issue6)
$ws.Activate() | Out-Null
